# Update res_bus vm_pu results for the 380 kV case (rows 2-25)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.063624978641885
$ws.Range("D2").Value = 1.065527254603892
$ws.Range("E2").Value = 1.07607684021187
$ws.Range("F2").Value = 1.080735208821614
$ws.Range("I2").Value = 1.049799515610171
$ws.Range("J2").Value = 1.06858939473689
$ws.Range("K2").Value = 1.06824025470723
$ws.Range("L2").Value = 1.078761724642886
$ws.Range("M2").Value = 1.083407866947428
$ws.Range("N2").Value = 1.02642921820604

# Row 3
$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.064959659140398
$ws.Range("D3").Value = 1.066556130166475
$ws.Range("E3").Value = 1.07734202337185
$ws.Range("F3").Value = 1.081947612596404
$ws.Range("I3").Value = 1.050146662359819
$ws.Range("J3").Value = 1.069577142982825
$ws.Range("K3").Value = 1.069084012965537
$ws.Range("L3").Value = 1.079843248133902
$ws.Range("M3").Value = 1.084437628448436
$ws.Range("N3").Value = 1.02676949508421

# Row 4
$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.06582281170965
$ws.Range("D4").Value = 1.067221220008263
$ws.Range("E4").Value = 1.078160483779806
$ws.Range("F4").Value = 1.082731777765999
$ws.Range("I4").Value = 1.050369407673085
$ws.Range("J4").Value = 1.070215305210042
$ws.Range("K4").Value = 1.069628711574627
$ws.Range("L4").Value = 1.08054230508067
$ws.Range("M4").Value = 1.085103036708704
$ws.Range("N4").Value = 1.026989080796131

# Row 5
$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.066185570740643
$ws.Range("D5").Value = 1.067500667223467
$ws.Range("E5").Value = 1.078504519752365
$ws.Range("F5").Value = 1.083061361431485
$ws.Range("I5").Value = 1.050462600113972
$ws.Range("J5").Value = 1.070483356942163
$ws.Range("K5").Value = 1.069857400727834
$ws.Range("L5").Value = 1.080836008718693
$ws.Range("M5").Value = 1.085382557088574
$ws.Range("N5").Value = 1.027081252564068

# Row 6
$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.066246473208191
$ws.Range("D6").Value = 1.067547578512172
$ws.Range("E6").Value = 1.078562282353746
$ws.Range("F6").Value = 1.083116695391404
$ws.Range("I6").Value = 1.050478221188888
$ws.Range("J6").Value = 1.070528350486498
$ws.Range("K6").Value = 1.069895780976062
$ws.Range("L6").Value = 1.080885312369
$ws.Range("M6").Value = 1.085429477111625
$ws.Range("N6").Value = 1.027096720297895

# Row 7
$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.065827659342783
$ws.Range("D7").Value = 1.067224954608555
$ws.Range("E7").Value = 1.078165080978108
$ws.Range("F7").Value = 1.082736181988089
$ws.Range("I7").Value = 1.05037065468099
$ws.Range("J7").Value = 1.070218887837957
$ws.Range("K7").Value = 1.069631768515324
$ws.Range("L7").Value = 1.080546230267094
$ws.Range("M7").Value = 1.085106772524741
$ws.Range("N7").Value = 1.026990312957201

# Row 8
$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.064076140556462
$ws.Range("D8").Value = 1.065875105959132
$ws.Range("E8").Value = 1.076504457254074
$ws.Range("F8").Value = 1.081145018376434
$ws.Range("I8").Value = 1.049917225629255
$ws.Range("J8").Value = 1.068923412496602
$ws.Range("K8").Value = 1.068525670659341
$ws.Range("L8").Value = 1.079127389542017
$ws.Range("M8").Value = 1.083756070663817
$ws.Range("N8").Value = 1.026544340212745

# Row 9
$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.060985950546801
$ws.Range("D9").Value = 1.06349135309598
$ws.Range("E9").Value = 1.073576597386515
$ws.Range("F9").Value = 1.078338468464863
$ws.Range("I9").Value = 1.049103775431842
$ws.Range("J9").Value = 1.066633037683048
$ws.Range("K9").Value = 1.066566783560792
$ws.Range("L9").Value = 1.076621278484578
$ws.Range("M9").Value = 1.081368848695684
$ws.Range("N9").Value = 1.025753883371209

# Row 10
$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.058923060311057
$ws.Range("D10").Value = 1.061898607656008
$ws.Range("E10").Value = 1.071623433004724
$ws.Range("F10").Value = 1.076465461333878
$ws.Range("I10").Value = 1.048551705616694
$ws.Range("J10").Value = 1.065100882152704
$ws.Range("K10").Value = 1.065254155829114
$ws.Range("L10").Value = 1.074946402007857
$ws.Range("M10").Value = 1.079772449807235
$ws.Range("N10").Value = 1.02522377945582

# Row 11
$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.058029102027826
$ws.Range("D11").Value = 1.061208056193058
$ws.Range("E11").Value = 1.070777355918633
$ws.Range("F11").Value = 1.075653926719289
$ws.Range("I11").Value = 1.04831032349836
$ws.Range("J11").Value = 1.064436167104268
$ws.Range("K11").Value = 1.064684157764381
$ws.Range("L11").Value = 1.074220147562011
$ws.Range("M11").Value = 1.079079994172885
$ws.Range("N11").Value = 1.024993485780704

# Row 12
$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.057696934694326
$ws.Range("D12").Value = 1.060951419689926
$ws.Range("E12").Value = 1.070463030245296
$ws.Range("F12").Value = 1.075352406978459
$ws.Range("I12").Value = 1.048220311880061
$ws.Range("J12").Value = 1.064189066928195
$ws.Range("K12").Value = 1.064472189147967
$ws.Range("L12").Value = 1.073950227815441
$ws.Range("M12").Value = 1.078822601609647
$ws.Range("N12").Value = 1.024907830112669

# Row 13
$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.057768190825899
$ws.Range("D13").Value = 1.061006475236276
$ws.Range("E13").Value = 1.070530456709788
$ws.Range("F13").Value = 1.07541708763922
$ws.Range("I13").Value = 1.048239635604893
$ws.Range("J13").Value = 1.064242079643678
$ws.Range("K13").Value = 1.064517668303316
$ws.Range("L13").Value = 1.074008133647769
$ws.Range("M13").Value = 1.078877821545977
$ws.Range("N13").Value = 1.024926208732774

# Row 14
$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.058001647270209
$ws.Range("D14").Value = 1.061186845306419
$ws.Range("E14").Value = 1.070751374800903
$ws.Range("F14").Value = 1.075629004645755
$ws.Range("I14").Value = 1.048302890288196
$ws.Range("J14").Value = 1.064415745729575
$ws.Range("K14").Value = 1.064666641403377
$ws.Range("L14").Value = 1.074197839121602
$ws.Range("M14").Value = 1.079058721786966
$ws.Range("N14").Value = 1.024986407794171

# Row 15
$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.058145472597729
$ws.Range("D15").Value = 1.061297959306106
$ws.Range("E15").Value = 1.070887482345127
$ws.Range("F15").Value = 1.075759563062879
$ws.Range("I15").Value = 1.048341816967405
$ws.Range("J15").Value = 1.064522721157023
$ws.Range("K15").Value = 1.064758395985573
$ws.Range("L15").Value = 1.074314702096653
$ws.Range("M15").Value = 1.079170155969713
$ws.Range("N15").Value = 1.025023483244552

# Row 16
$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.058982373987499
$ws.Range("D16").Value = 1.06194441849088
$ws.Range("E16").Value = 1.07167957682617
$ws.Range("F16").Value = 1.076519309208823
$ws.Range("I16").Value = 1.048567676116481
$ws.Range("J16").Value = 1.065144969855311
$ws.Range("K16").Value = 1.065291950359751
$ws.Range("L16").Value = 1.074994579339478
$ws.Range("M16").Value = 1.079818380196611
$ws.Range("N16").Value = 1.025239047302275

# Row 17
$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.059507146163167
$ws.Range("D17").Value = 1.062349687380377
$ws.Range("E17").Value = 1.072176342347614
$ws.Range("F17").Value = 1.076995739191326
$ws.Range("I17").Value = 1.048708726369622
$ws.Range("J17").Value = 1.065534944886162
$ws.Range("K17").Value = 1.065626199134539
$ws.Range("L17").Value = 1.075420772684063
$ws.Range("M17").Value = 1.080224669734004
$ws.Range("N17").Value = 1.025374062132903

# Row 18
$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.059813168372133
$ws.Range("D18").Value = 1.06258598891415
$ws.Range("E18").Value = 1.072466064302762
$ws.Range("F18").Value = 1.077273583720834
$ws.Range("I18").Value = 1.048790773647273
$ws.Range("J18").Value = 1.065762287263922
$ws.Range("K18").Value = 1.065821004488446
$ws.Range("L18").Value = 1.075669265658171
$ws.Range("M18").Value = 1.08046153553055
$ws.Range("N18").Value = 1.025452741195721

# Row 19
$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.059917502543023
$ws.Range("D19").Value = 1.062666547293222
$ws.Range("E19").Value = 1.072564846509005
$ws.Range("F19").Value = 1.077368313301
$ws.Range("I19").Value = 1.048818711528179
$ws.Range("J19").Value = 1.065839784284839
$ws.Range("K19").Value = 1.06588740163838
$ws.Range("L19").Value = 1.075753978741115
$ws.Range("M19").Value = 1.080542281050539
$ws.Range("N19").Value = 1.025479556400709

# Row 20
$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.05945085019508
$ws.Range("D20").Value = 1.062306214669501
$ws.Range("E20").Value = 1.072123047543688
$ws.Range("F20").Value = 1.076944627856291
$ws.Range("I20").Value = 1.04869361628875
$ws.Range("J20").Value = 1.065493117047901
$ws.Range("K20").Value = 1.065590353589108
$ws.Range("L20").Value = 1.075375056371651
$ws.Range("M20").Value = 1.080181090732997
$ws.Range("N20").Value = 1.025359583856562

# Row 21
$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.057932903276356
$ws.Range("D21").Value = 1.061133734560853
$ws.Range("E21").Value = 1.070686321473285
$ws.Range("F21").Value = 1.075566602592205
$ws.Range("I21").Value = 1.048284273074443
$ws.Range("J21").Value = 1.064364610813168
$ws.Range("K21").Value = 1.064622779346833
$ws.Range("L21").Value = 1.074141979932289
$ws.Range("M21").Value = 1.079005456268362
$ws.Range("N21").Value = 1.024968683834582

# Row 22
$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.056977863009681
$ws.Range("D22").Value = 1.060395767839757
$ws.Range("E22").Value = 1.069782675078682
$ws.Range("F22").Value = 1.074699721076656
$ws.Range("I22").Value = 1.048024868197976
$ws.Range("J22").Value = 1.063653942137919
$ws.Range("K22").Value = 1.064013002828604
$ws.Range("L22").Value = 1.073365788833646
$ws.Range("M22").Value = 1.078265223955729
$ws.Range("N22").Value = 1.024722247479185

# Row 23
$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.057484210664048
$ws.Range("D23").Value = 1.060787052887062
$ws.Range("E23").Value = 1.070261746533988
$ws.Range("F23").Value = 1.07515931620109
$ws.Range("I23").Value = 1.048162576871511
$ws.Range("J23").Value = 1.064030789310813
$ws.Range("K23").Value = 1.064336392675777
$ws.Range("L23").Value = 1.073777349393523
$ws.Range("M23").Value = 1.078657736916524
$ws.Range("N23").Value = 1.024852951100231

# Row 24
$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.059476288131378
$ws.Range("D24").Value = 1.062325858378754
$ws.Range("E24").Value = 1.072147129273328
$ws.Range("F24").Value = 1.076967723019939
$ws.Range("I24").Value = 1.048700444578926
$ws.Range("J24").Value = 1.065512017628941
$ws.Range("K24").Value = 1.065606551131031
$ws.Range("L24").Value = 1.075395713909337
$ws.Range("M24").Value = 1.080200782566128
$ws.Range("N24").Value = 1.025366126191374

# Row 25
$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.061785310000376
$ws.Range("D25").Value = 1.064108232223333
$ws.Range("E25").Value = 1.074333729057744
$ws.Range("F25").Value = 1.079064365927491
$ws.Range("I25").Value = 1.049315789784432
$ws.Range("J25").Value = 1.067226067466837
$ws.Range("K25").Value = 1.067074375942775
$ws.Range("L25").Value = 1.077269886388566
$ws.Range("M25").Value = 1.081986859579553
$ws.Range("N25").Value = 1.025958784427353
